# ocd convert folder names to consistent lowercase
# Walk every slide/shape and, for any picture whose AlternativeText
# (maps to the OOXML cNvPr "descr" attribute) references the
# "../Images/..." folder, rewrite it to use the lowercase "../images/..."
# folder so all image references are consistently cased.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        $alt = $sh.AlternativeText
        if ($alt -and $alt.Contains("../Images/")) {
            $sh.AlternativeText = $alt.Replace("../Images/", "../images/")
        }
    }
}
